$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that needs to be
# bumped from 46075 to 46076 for every data row (rows 2 through 95).
for ($r = 2; $r -le 95; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
